# Horario informática — update for the week of "6 hasta el 10 de marzo" schedule block
# (rows 69-88 on Hoja1): lunes was spent entirely advancing the lab practice, martes followed
# the plan, miercoles the person overslept, jueves had to rush to finish the lab practice
# (due friday 10 de marzo), and friday followed the plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122  # xlPasteFormats

# --- Stable reference cells whose style never changes, used as format donors ---
# style index 1 (border only, no fill)      -> C72
# style index 4 (border + dark/black fill)  -> B70
# style index 6 (border + light/white fill) -> F82

function Set-Style($cellRef, $sourceRef) {
    $ws.Range($sourceRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------- Style-only changes first (values unaffected by these) ----------
Set-Style "F76" "B70"   # s6 -> s4
Set-Style "C78" "B70"   # s1 -> s4
Set-Style "D78" "F82"   # s1 -> s6
Set-Style "E78" "B70"   # s1 -> s4
Set-Style "F78" "B70"   # s1 -> s4
Set-Style "C79" "B70"   # s1 -> s4
Set-Style "D79" "B70"   # s1 -> s4
Set-Style "D80" "B70"   # s1 -> s4
Set-Style "D81" "B70"   # s6 -> s4

# ---------- Style + value changes ----------
Set-Style "B75" "F82"   # s4 -> s6
$ws.Range("B75").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "B76" "F82"   # s4 -> s6
$ws.Range("B76").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "E76" "F82"   # (no fill style) -> s6
$ws.Range("E76").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "B77" "F82"   # (no fill style) -> s6
$ws.Range("B77").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "E77" "F82"   # (no fill style) -> s6
$ws.Range("E77").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "B81" "F82"   # s4 -> s6
$ws.Range("B81").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "F81" "F82"   # (no fill style) -> s6
$ws.Range("F81").Value = "estudiar/prácticas/tareas informática"

Set-Style "B82" "F82"   # s4 -> s6
$ws.Range("B82").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "D82" "F82"   # s4 -> s6
$ws.Range("D82").Value = "estudiar algebra lineal"

Set-Style "B83" "F82"   # (no fill style) -> s6
$ws.Range("B83").Value = "Adelantar prácticas laboratorio informatica"

Set-Style "D83" "F82"   # s4 -> s6
$ws.Range("D83").Value = "estudiar algebra lineal"

Set-Style "B86" "B70"   # s1 -> s4
# value of B86 (14, "estudiar/p+racticas/tareas informática") is unchanged

Set-Style "E88" "F82"   # (no fill style) -> s6
$ws.Range("E88").Value = "estudiar cálculo integral"

# ---------- Value-only changes (style stays the same) ----------
$ws.Range("B79").Value = "Adelantar prácticas laboratorio informatica"
$ws.Range("E79").Value = "Adelantar prácticas laboratorio informatica"
$ws.Range("F79").Value = "estudiar/prácticas/tareas informática"

$ws.Range("B80").Value = "Adelantar prácticas laboratorio informatica"
$ws.Range("E80").Value = "clase de física mecánica"

$ws.Range("C86").Value = "estudiar algebrea lineal"
$ws.Range("E86").Value = "Adelantar prácticas laboratorio informatica"

$ws.Range("C87").Value = "estudiar algebra lineal"
$ws.Range("E87").Value = "Adelantar prácticas laboratorio informatica"

# ---------- Update the saved view/selection ----------
# (topLeftCell reflects the window's scroll position; activeCell/sqref reflect the selection)
$win = $excel.ActiveWindow
$win.ScrollRow = 70
$win.ScrollColumn = 1
$ws.Range("F90").Select() | Out-Null
